# feat: activity to read, write and append from excel
#
# Rewrites the small Salary/Year demo table into a Years/Salary table
# (columns swapped, header text tweaked, values refreshed) and appends
# one more row of data, mimicking a "read, write, append" RPA exercise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: column order flips to Years | Salary, and the "Year" label
# becomes "Years".
$ws.Range("A1").Value = "Years"
$ws.Range("B1").Value = "Salary"

# Full Years/Salary data set (5 existing rows refreshed + 1 appended row).
$data = @(
    @(2019, 30000),
    @(2020, 30000),
    @(2021, 35000),
    @(2022, 40000),
    @(2023, 50000)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Leave the selection where it would land right after appending the last
# row (one row below the newly written data, in the Salary column).
$lastRow = $startRow + $data.Length
$ws.Range("B$lastRow").Select()
